$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add new header values in columns P and Q,
#     copying the formatting (style) already used by the O1 header cell ---
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14

$ws.Range("O1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Value = 15

# --- Rows 2-25: fix columns I, K, M, O and add columns P, Q ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I column
    $ws.Cells.Item($r, 11).Value = 1   # K column
    $ws.Cells.Item($r, 13).Value = 2   # M column
    $ws.Cells.Item($r, 15).Value = 1   # O column
    $ws.Cells.Item($r, 16).Value = 2   # P column (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q column (new)
}
